# Apply the xlsx-handlebars "deleteCurrentSheet / setCurrentSheetName /
# hideCurrentSheet / hyperlink" example edits to the workbook.
# Shared-string order matters for matching the canonical OOXML, so cell
# values below are written in the same order the original author's Excel
# session would have produced them.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. "被删除的工作表" (deleteCurrentSheet demo sheet)
# ---------------------------------------------------------------------
$sheetDeleted = $wb.Worksheets.Add($null, $ws1)
$sheetDeleted.Name = "被删除的工作表"
$sheetDeleted.Range("A1").Value = "{{deleteCurrentSheet}}"
$sheetDeleted.Columns("A").ColumnWidth = 18.0375

# ---------------------------------------------------------------------
# 2. setCurrentSheetName demo row on Sheet1
# ---------------------------------------------------------------------
$ws1.Range("A24").Value = "工作表命名:"
$ws1.Range("B24").Value = '{{setCurrentSheetName "new sheet name"}}'

# ---------------------------------------------------------------------
# 3. "被隐藏的工作表" (hideCurrentSheet demo sheet)
# ---------------------------------------------------------------------
$sheetHidden = $wb.Worksheets.Add($null, $sheetDeleted)
$sheetHidden.Name = "被隐藏的工作表"
$sheetHidden.Range("A1").Value = "{{hideCurrentSheet}}"
$sheetHidden.Columns("A").ColumnWidth = 18.0375

# ---------------------------------------------------------------------
# 4. "被链接的工作表" + hyperlink from Sheet1!A26
# ---------------------------------------------------------------------
$sheetLinked = $wb.Worksheets.Add($null, $sheetHidden)
$sheetLinked.Name = "被链接的工作表"

$ws1.Hyperlinks.Add($ws1.Range("A26"), "", "被链接的工作表!A1", "", "链接")

# ---------------------------------------------------------------------
# 5. Static (non-templated) merged header cell, H1:K1
# ---------------------------------------------------------------------
$ws1.Range("H1:K1").Merge()
$ws1.Range("H1:K1").HorizontalAlignment = -4108
$ws1.Range("H1").Value = "静态合并单元格"

# ---------------------------------------------------------------------
# 6. Selection / active cell, matching the edited template
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("H1:K1").Select()
